# Commit: Update 11/10/2022 BFF-Users GET Method
#
# 1) Users sheet: add CUIT 1-4 / Rol / CID CUIT 1-4 columns (C:K) to the
#    header row and to the existing "ADMIN" detail row (row 19).
$wb = $excel.ActiveWorkbook
$users = $wb.Worksheets.Item("Users")

$users.Range("C1").Value = "CUIT 1"
$users.Range("D1").Value = "CUIT 2"
$users.Range("E1").Value = "CUIT 3"
$users.Range("F1").Value = "CUIT 4"
$users.Range("G1").Value = "Rol"
$users.Range("H1").Value = "CID CUIT 1"
$users.Range("I1").Value = "CID CUIT 2"
$users.Range("J1").Value = "CID CUIT 3"
$users.Range("K1").Value = "CID CUIT 4"

$users.Range("A19").ClearFormats()
$users.Range("C19").Value = "30-61398988-5"
$users.Range("D19").Value = "30-53725008-5"
$users.Range("E19").Value = "20-10474244-1"
$users.Range("F19").Value = "20-07608479-4"
$users.Range("G19").Value = "ADMIN"
$users.Range("H19").Value = 1934
$users.Range("I19").Value = 1933
$users.Range("J19").Value = 1288
$users.Range("K19").Value = 1809

$users.Range("C1:F19").ColumnWidth = 13.42578125
$users.Range("H1").ColumnWidth = 12.85546875

$users.Range("A19").Select() | Out-Null

# 2) Insert a new "Roles" worksheet right after "Users" describing each
#    Keycloak role exposed to the BFF.
$roles = $wb.Worksheets.Add($null, $users)
$roles.Name = "Roles"
$roles.Range("A1").Value = "Nombre"
$roles.Range("B1").Value = "Id"
$roles.Range("A2").Value = "Acceso offline"
$roles.Range("B2").Value = "ROLE_offline_access"
$roles.Range("A3").Value = "UMA Authorization"
$roles.Range("B3").Value = "ROLE_uma_authorization"
$roles.Range("A4").Value = "Admin"
$roles.Range("B4").Value = "ROLE_admin"
$roles.Range("A5").Value = "PayWay Free"
$roles.Range("B5").Value = "ROLE_phe-free"
$roles.Range("A1").ColumnWidth = 18
$roles.Range("B9").Select() | Out-Null

# 3) EndPoints sheet: the former "PHE/Security-certs-controller-adapter/..."
#    entry is replaced by the BFF "Users/Users" endpoint, and a new
#    "Users/BFF-User" (BFF User) endpoint row is appended.
$endpoints = $wb.Worksheets.Item("EndPoints")
$endpoints.Range("B3").Value = "Users/Users"
$endpoints.Range("A4").Value = "BFF User"
$endpoints.Range("B4").Value = "Users/BFF-User"
$endpoints.Range("A8").Select() | Out-Null

# 4) HTTPCodes sheet: add the 409 Conflict status code.
$httpcodes = $wb.Worksheets.Item("HTTPCodes")
$httpcodes.Range("A5").Value = "Conflict"
$httpcodes.Range("B5").Value = 409
$httpcodes.Range("C10").Select() | Out-Null

# 5) Servers sheet: only the remembered selection moved (no data changed).
$servers = $wb.Worksheets.Item("Servers")
$servers.Range("B2").Select() | Out-Null

# 6) Make "Roles" the active tab, matching the saved workbook view.
$roles.Activate() | Out-Null
